$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9441280961036682
$ws.Range("B1").Value = 1.952826738357544
$ws.Range("C1").Value = 4.278550624847412
$ws.Range("D1").Value = 3.295520067214966
$ws.Range("E1").Value = 1.443691492080688
